# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") holds strikeouts for each outing. The prior export derived
# this column from a different (incorrect) source; this script rewrites it
# with the correct K values for every game row (rows 2-66) on the active
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ("K") values, in order, for rows 2 through 66.
$s_vals = @(
    2, 2, 1, 4, 3, 1, 0, 3, 1, 2,
    0, 0, 1, 2, 2, 2, 0, 0, 3, 2,
    0, 2, 2, 1, 2, 2, 1, 4, 1, 2,
    1, 4, 1, 0, 1, 1, 1, 1, 0, 1,
    1, 2, 1, 2, 1, 0, 1, 1, 0, 1,
    0, 1, 4, 3, 0, 1, 1, 3, 2, 1,
    5, 2, 3, 4, 3
)

$startRow = 2
$col = 7  # column G = K

for ($i = 0; $i -lt $s_vals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, $col).Value = $s_vals[$i]
}
